$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, row 1 (next to the existing "sum" header in G1)
$ws.Range("H1").Value = "Save"

# Copy the header formatting (bold, centered, bordered) from the neighboring G1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the corresponding data value for row 2 in the new column
$ws.Range("H2").Value = 0
